$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'61.390.24"
$ws.Cells.Item(2, 5).Value = "  +0.46%  "

$ws.Cells.Item(3, 4).Value = "'2.931.18"
$ws.Cells.Item(3, 5).Value = "  +0.39%  "

$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.09%  "

$ws.Cells.Item(5, 4).Value = "'595.04"
$ws.Cells.Item(5, 5).Value = "  +0.93%  "

$ws.Cells.Item(6, 4).Value = "'144.97"
$ws.Cells.Item(6, 5).Value = "  -0.07%  "

$ws.Cells.Item(7, 5).Value = "  -0.04%  "

$ws.Cells.Item(8, 5).Value = "  -0.38%  "

$ws.Cells.Item(9, 4).Value = "'6.95"
$ws.Cells.Item(9, 5).Value = "  +2.44%  "

$ws.Cells.Item(10, 4).Value = "'0.142"
$ws.Cells.Item(10, 5).Value = "  -1.18%  "

$ws.Cells.Item(11, 4).Value = "'0.441"
$ws.Cells.Item(11, 5).Value = "  -0.32%  "

$ws.Cells.Item(12, 5).Value = "  -0.48%  "

$ws.Cells.Item(13, 4).Value = "'33.70"
$ws.Cells.Item(13, 5).Value = "  -0.04%  "

$ws.Cells.Item(14, 5).Value = "  +0.73%  "

$ws.Cells.Item(15, 4).Value = "'3.414.74"
$ws.Cells.Item(15, 5).Value = "  +0.26%  "

$ws.Cells.Item(16, 4).Value = "'61.335.23"
$ws.Cells.Item(16, 5).Value = "  +0.45%  "

$ws.Cells.Item(17, 5).Value = "  +0.40%  "

$ws.Cells.Item(18, 4).Value = "'2.929.24"
$ws.Cells.Item(18, 5).Value = "  +0.23%  "

$ws.Cells.Item(19, 4).Value = "'432.62"
$ws.Cells.Item(19, 5).Value = "  +0.40%  "

$ws.Cells.Item(20, 4).Value = "'13.52"
$ws.Cells.Item(20, 5).Value = "  +0.49%  "

$ws.Cells.Item(21, 5).Value = "  -0.28%  "

$ws.Cells.Item(22, 5).Value = "  +0.90%  "

$ws.Cells.Item(23, 4).Value = "'81.83"
$ws.Cells.Item(23, 5).Value = "  +1.09%  "

$ws.Cells.Item(24, 4).Value = "'10.88"
$ws.Cells.Item(24, 5).Value = "  -0.75%  "

$ws.Cells.Item(25, 4).Value = "'2.20"
$ws.Cells.Item(25, 5).Value = "  -1.13%  "

$ws.Cells.Item(26, 5).Value = "  -1.74%  "

$ws.Cells.Item(27, 4).Value = "'1.00"
$ws.Cells.Item(27, 5).Value = "  -0.07%  "

$ws.Cells.Item(28, 5).Value = "  -2.74%  "

$ws.Cells.Item(29, 5).Value = "  -0.25%  "

$ws.Cells.Item(30, 4).Value = "'6.94"
$ws.Cells.Item(30, 5).Value = "  -2.64%  "

$ws.Cells.Item(31, 5).Value = "  +1.69%  "

$ws.Cells.Item(32, 4).Value = "'26.66"
$ws.Cells.Item(32, 5).Value = "  +0.58%  "

$ws.Cells.Item(33, 5).Value = "  -0.02%  "

$ws.Cells.Item(34, 4).Value = "'0.0₃0881"
$ws.Cells.Item(34, 5).Value = "  +2.11%  "

$ws.Cells.Item(35, 5).Value = "  +0.25%  "

$ws.Cells.Item(36, 5).Value = "  +0.28%  "

$ws.Cells.Item(37, 5).Value = "  -2.01%  "

$ws.Cells.Item(38, 2).Value = "Stacks"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(38, 4).Value = "'2.01"
$ws.Cells.Item(38, 5).Value = "  +0.47%  "

$ws.Cells.Item(39, 2).Value = "Kaspa"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(39, 4).Value = "'0.124"
$ws.Cells.Item(39, 5).Value = "  -0.63%  "

$ws.Cells.Item(40, 4).Value = "'8.63"
$ws.Cells.Item(40, 5).Value = "  +0.36%  "

$ws.Cells.Item(41, 5).Value = "  +8.42%  "

$ws.Cells.Item(42, 5).Value = "  -1.57%  "

$ws.Cells.Item(43, 4).Value = "'0.0348"
$ws.Cells.Item(43, 5).Value = "  +0.19%  "

$ws.Cells.Item(44, 4).Value = "'2.707.22"
$ws.Cells.Item(44, 5).Value = "  -0.10%  "

$ws.Cells.Item(45, 4).Value = "'368.05"
$ws.Cells.Item(45, 5).Value = "  -2.41%  "

$ws.Cells.Item(46, 4).Value = "'134.36"
$ws.Cells.Item(46, 5).Value = "  +2.03%  "

$ws.Cells.Item(48, 4).Value = "'23.76"
$ws.Cells.Item(48, 5).Value = "  -1.67%  "

$ws.Cells.Item(49, 5).Value = "  -1.20%  "

$ws.Cells.Item(50, 5).Value = "  -1.36%  "

$ws.Cells.Item(51, 5).Value = "  -0.57%  "
